$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z1").Formula = "=""63.664.44"""
$ws.Range("Z1").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("Z1").Formula = "=""3.403.21"""
$ws.Range("Z1").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("Z1").Formula = "=""568.06"""
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("Z1").Formula = "=""155.81"""
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E6").Value = "  +2.11%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("Z1").Formula = "=""3.403.99"""
$ws.Range("Z1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("Z1").Formula = "=""0.575"""
$ws.Range("Z1").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E9").Value = "  +8.77%  "

$ws.Range("Z1").Formula = "=""7.33"""
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E10").Value = "  -1.18%  "

$ws.Range("E11").Value = "  +4.22%  "

$ws.Range("Z1").Formula = "=""0.441"""
$ws.Range("Z1").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("Z1").Formula = "=""3.991.24"""
$ws.Range("Z1").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E13").Value = "  +1.36%  "

$ws.Range("E14").Value = "  -3.15%  "

$ws.Range("Z1").Formula = "=""0.0000190"""
$ws.Range("Z1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E15").Value = "  +5.61%  "

$ws.Range("Z1").Formula = "=""27.46"""
$ws.Range("Z1").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("Z1").Formula = "=""63.682.95"""
$ws.Range("Z1").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("Z1").Formula = "=""3.384.66"""
$ws.Range("Z1").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("Z1").Formula = "=""6.35"""
$ws.Range("Z1").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("Z1").Formula = "=""14.15"""
$ws.Range("Z1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E20").Value = "  +1.69%  "

$ws.Range("Z1").Formula = "=""383.04"""
$ws.Range("Z1").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("Z1").Formula = "=""8.08"""
$ws.Range("Z1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E22").Value = "  -3.90%  "

$ws.Range("Z1").Formula = "=""73.14"""
$ws.Range("Z1").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E23").Value = "  +3.66%  "

$ws.Range("Z1").Formula = "=""0.998"""
$ws.Range("Z1").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("Z1").Formula = "=""0.0000119"""
$ws.Range("Z1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E26").Value = "  +23.14%  "

$ws.Range("Z1").Formula = "=""9.53"""
$ws.Range("Z1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E27").Value = "  +1.97%  "

$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("Z1").Formula = "=""6.09"""
$ws.Range("Z1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E30").Value = "  +9.75%  "

$ws.Range("Z1").Formula = "=""1.40"""
$ws.Range("Z1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E31").Value = "  +7.13%  "

$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("Z1").Formula = "=""23.34"""
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("Z1").Formula = "=""6.41"""
$ws.Range("Z1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E34").Value = "  +2.54%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("Z1").Formula = "=""6.88"""
$ws.Range("Z1").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E36").Value = "  +2.86%  "

$ws.Range("Z1").Formula = "=""159.85"""
$ws.Range("Z1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E37").Value = "  +1.30%  "

$ws.Range("E38").Value = "  -0.41%  "

$ws.Range("Z1").Formula = "=""0.0762"""
$ws.Range("Z1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("Z1").Formula = "=""1.84"""
$ws.Range("Z1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E40").Value = "  -1.57%  "

$ws.Range("Z1").Formula = "=""2.891.94"""
$ws.Range("Z1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E41").Value = "  +1.20%  "

$ws.Range("Z1").Formula = "=""26.76"""
$ws.Range("Z1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E42").Value = "  -0.57%  "

$ws.Range("Z1").Formula = "=""0.0316"""
$ws.Range("Z1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E43").Value = "  -3.20%  "

$ws.Range("Z1").Formula = "=""42.27"""
$ws.Range("Z1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E44").Value = "  +3.26%  "

$ws.Range("E45").Value = "  +2.58%  "

$ws.Range("Z1").Formula = "=""0.755"""
$ws.Range("Z1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("Z1").Formula = "=""23.21"""
$ws.Range("Z1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E47").Value = "  +5.81%  "

$ws.Range("E48").Value = "  +2.65%  "

$ws.Range("Z1").Formula = "=""2.19"""
$ws.Range("Z1").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E49").Value = "  +21.23%  "

$ws.Range("E50").Value = "  +3.52%  "

$ws.Range("Z1").Formula = "=""6.44"""
$ws.Range("Z1").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E51").Value = "  +2.22%  "

$excel.CutCopyMode = 0
